$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H13").Value = 907.5
$ws1.Range("I13").Value = 377.14
$ws1.Range("H29").Value = "1 de 27"
$ws1.Range("I29").Value = "1 de 27"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 1284.64
$ws2.Range("F29").Value = 5825.46

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D7").Value = 907.5
$ws3.Range("E7").Value = -347.5
$ws3.Range("F7").Value = 1.620535714285714

$ws3.Range("D8").Value = 377.14
$ws3.Range("E8").Value = 247.86
$ws3.Range("F8").Value = 0.603424

$ws3.Range("D19").Value = 5825.46
$ws3.Range("E19").Value = 31674.54093005039
$ws3.Range("F19").Value = 0.1553455961472205
